$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the homework 1 (column C) score to 5 for the specified students/rows
$rows = @(10, 11, 15, 16, 19, 20, 21, 23, 25)
foreach ($r in $rows) {
    $ws.Range("C$r").Value = 5
}

# Update the active selection in the frozen (bottom-right) pane to C12
$ws.Range("C12").Select()
